# Update dSF (column F) values on Sheet1 to reflect repulled data / recalculated mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F6"  = 0
    "F15" = -3
    "F17" = -2
    "F22" = -5
    "F24" = -5
    "F26" = 0
    "F27" = 2
    "F29" = -6
    "F31" = 10
    "F33" = 8
    "F34" = 6
    "F37" = -3
    "F41" = 0
    "F43" = 1
    "F44" = 1
    "F46" = -2
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

$wb.Save()
